$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.388.98"
$ws.Range("E2").Value = "  -3.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.462.83"
$ws.Range("E3").Value = "  -5.48%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.79"
$ws.Range("E5").Value = "  -4.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.61"
$ws.Range("E6").Value = "  -6.60%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  -5.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.459.52"
$ws.Range("E9").Value = "  -5.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.104"
$ws.Range("E10").Value = "  -11.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.153"
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.26"
$ws.Range("E12").Value = "  -9.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("E13").Value = "  -9.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.32"
$ws.Range("E14").Value = "  -10.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.911.16"
$ws.Range("E15").Value = "  -5.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.322.15"
$ws.Range("E16").Value = "  -3.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000159"
$ws.Range("E17").Value = "  -10.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.489.16"
$ws.Range("E18").Value = "  -3.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.98"
$ws.Range("E19").Value = "  -8.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.83"
$ws.Range("E20").Value = "  -9.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.10"
$ws.Range("E21").Value = "  -9.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "314.60"
$ws.Range("E22").Value = "  -8.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.92"
$ws.Range("E24").Value = "  -6.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.70"
$ws.Range("E25").Value = "  -4.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000100"
$ws.Range("E26").Value = "  -7.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.623.46"
$ws.Range("E27").Value = "  -3.93%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "527.79"
$ws.Range("E29").Value = "  -10.24%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.21"
$ws.Range("E30").Value = "  -9.89%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("E31").Value = "  -7.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.39"
$ws.Range("E32").Value = "  -6.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.143"
$ws.Range("E33").Value = "  -10.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.87"
$ws.Range("E34").Value = "  -9.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  -10.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.70"
$ws.Range("E36").Value = "  -12.55%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.997"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.69"
$ws.Range("E38").Value = "  -13.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.371"
$ws.Range("E39").Value = "  -8.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.15"
$ws.Range("E40").Value = "  -7.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "140.88"
$ws.Range("E41").Value = "  -9.08%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.65"
$ws.Range("E43").Value = "  -11.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.54"
$ws.Range("E44").Value = "  -1.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.22"
$ws.Range("E45").Value = "  -8.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "145.65"
$ws.Range("E46").Value = "  -6.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.50"
$ws.Range("E47").Value = "  -10.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.66"
$ws.Range("E48").Value = "  -10.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0527"
$ws.Range("E49").Value = "  -10.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.581"
$ws.Range("E50").Value = "  -7.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0930"
$ws.Range("E51").Value = "  -7.20%  "
